# [ResourceUnit/Table] Add a "UnitName" (FString) column to the resource
# unit table, insert a type-declaration row under the headers, and fill
# in sample duck names for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The trailing, never-used F/G columns on row 1 are not needed any more
# once the sheet only spans to column E - drop them outright.
$ws.Range("F1:G1").Clear()

# --- Make room for the new type-declaration row under the header row
#     (row 4). A plain row insert only moves rows, so row 1/2 above and
#     column layout are unaffected. ---
$ws.Range("4:4").Insert()

# Row 4 (now the blank, freshly-inserted row) used to be the first data
# row ("1,1,-1"); that data lives on row 5 now. Move the
# "ChildTableDataType" column's old values (which were in column C) over
# to their new home in column D, then drop in the new UnitName values in
# column C, freeing column C up for the new field everywhere it's used.
$ws.Range("D3").Value = "ChildTableDataType"
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 12

# --- Column header (row 3) ---
$ws.Range("C3").Value = "UnitName"

# --- Type-declaration row (row 4) ---
$ws.Range("A4").Value = "int32"
$ws.Range("B4").Value = "int32"
$ws.Range("C4").Value = "FString"
$ws.Range("D4").Value = "int32"

# --- Data rows: new UnitName values ---
$ws.Range("C5").Value = "흰오리"
$ws.Range("C6").Value = "노오란오리"
$ws.Range("C7").Value = "푸르스름오리"

# --- Widen the new column to fit its content ---
$ws.Range("C:C").Columns.AutoFit()

# --- Style row 2 (the ChildTableDataType description row) with the
#     filled "메모" (Note) look: font on A2, fill+border across A2:E2 ---
$ws.Range("A2:E2").Interior.Color = 13434879
$ws.Range("A2:E2").Borders.LineStyle = 1
$ws.Range("A2:E2").Borders.Color = 11776946

# --- Page setup + selection to mirror the authored workbook ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("C8").Select()
